$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column M (year 2021) to the table ---

# M3: blank cell, same formatting as the other cells in the thick-bottom-border row (same as L3)
$ws.Range("L3").Copy()
$ws.Range("M3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# M4: header year value, bold Times New Roman 9pt with a medium bottom border,
# vertically bottom-aligned (no explicit center alignment)
$ws.Cells.Item(4, 13).Value = 2021
$m4 = $ws.Range("M4")
$m4.ClearFormats()
$m4.Font.Name = "Times New Roman"
$m4.Font.Family = 1
$m4.Font.Size = 9
$m4.Font.Bold = $true
$m4.Font.ThemeColor = 1
$m4.Borders.Item(9).Weight = -4138   # xlMedium bottom border

# M5 / M6: data values, regular Times New Roman 9pt, no border, no special alignment
$ws.Cells.Item(5, 13).Value = 93.5
$m5 = $ws.Range("M5")
$m5.ClearFormats()
$m5.Font.Name = "Times New Roman"
$m5.Font.Family = 1
$m5.Font.Size = 9
$m5.Font.Bold = $false
$m5.Font.ThemeColor = 1

$ws.Cells.Item(6, 13).Value = 96.6
$m6 = $ws.Range("M6")
$m6.ClearFormats()
$m6.Font.Name = "Times New Roman"
$m6.Font.Family = 1
$m6.Font.Size = 9
$m6.Font.Bold = $false
$m6.Font.ThemeColor = 1

# M7: data value, regular Times New Roman 9pt with a medium bottom border
$ws.Cells.Item(7, 13).Value = 98.1
$m7 = $ws.Range("M7")
$m7.ClearFormats()
$m7.Font.Name = "Times New Roman"
$m7.Font.Family = 1
$m7.Font.Size = 9
$m7.Font.Bold = $false
$m7.Font.ThemeColor = 1
$m7.Borders.Item(9).Weight = -4138   # xlMedium bottom border

# --- Update the view: scroll back to A1 and move the selection ---
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("N11").Select()
